# "Hero object it is passed to backstory methods"
#
# The "Przeszłość" (backstory) sheet gets a new column C holding the
# mechanical "addition" (effect) that corresponds to each backstory
# row in column B - this is what a Hero object consumes when a
# backstory method is applied. A header label is added in C1, and the
# rows that grant a concrete bonus (Splugawienie/Szaleństwo points,
# extra children, an extra language) get their effect spelled out.

$wb = $excel.ActiveWorkbook

$wsBackstory = $wb.Worksheets.Item("Przeszłość")

$wsBackstory.Range("C3").Value = "1 punkt splugawienia"
$wsBackstory.Range("C5").Value = "1 punkt splugawienia"
$wsBackstory.Range("C7").Value = "1 punkt szaleństwa"
$wsBackstory.Range("C13").Value = "1k6-2 dzieci"
$wsBackstory.Range("C14").Value = "Mówisz w dodatkowym języku"
$wsBackstory.Range("C1").Value = "addition"

# Selections left behind by the editing session.
$wsHumanAncestry = $wb.Worksheets.Item("humanAncestry")
$wsHumanAncestry.Activate() | Out-Null
$wsHumanAncestry.Range("F21").Select() | Out-Null

$wsBackstory.Activate() | Out-Null
$wsBackstory.Range("D19").Select() | Out-Null

$wsAge = $wb.Worksheets.Item("Wiek")
$wsAge.Activate() | Out-Null
$wsAge.Range("B9").Select() | Out-Null

# Restore the tab that was active/visible when the workbook was saved.
$wsHumanAncestry.Activate() | Out-Null
